$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: transaction-count ratio (column H) "2:0" -> "1:0"
$ws.Range("H7").Value = "1:0"

# Row 8: transaction-count ratio (H) "2:0" -> "1:1";
# sale price (P) "50.0000" -> "25.0000" (stored as text, so temporarily force
# a text number format to avoid Excel auto-converting the literal to a
# number, then restore the original numeric format so the cell style is
# unchanged);
# transactions ratio (Q) "1:0" -> "0:1"
$ws.Range("H8").Value = "1:1"

$origFormatP8 = $ws.Range("P8").NumberFormat
$ws.Range("P8").NumberFormat = "@"
$ws.Range("P8").Value = "25.0000"
$ws.Range("P8").NumberFormat = $origFormatP8

$ws.Range("Q8").Value = "0:1"

# Row 9: transaction-count ratio (H) "1:0" -> "0:0"
$ws.Range("H9").Value = "0:0"

# Row 10: transaction-count ratio (H) "2:0" -> "1:0"
$ws.Range("H10").Value = "1:0"

# Row 11: transaction-count ratio (H) "2:0" -> "1:0"
$ws.Range("H11").Value = "1:0"

# Row 12: total (N) recalculated from 1154 to 1129 (reflects the P8 price drop)
$ws.Range("N12").Value = 1129
